# Commit: "Add files via upload"
# Adds a new worksheet "Gonsalves 204" (quotation #204) at the end of the
# workbook, containing a two-option (Honeywell / Dahua) CCTV quotation
# table, and makes it the active/selected sheet - matching the author's
# edit of uploading one more quotation sheet to the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new worksheet as the LAST tab in the workbook
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Gonsalves 204"

# Column widths matching the authored sheet
$ws.Columns.Item(2).ColumnWidth = 38.21875
$ws.Columns.Item(3).ColumnWidth = 13.109375

# ---------------------------------------------------------------------
# Helper-ish inline style setup (re-used for both option tables)
# ---------------------------------------------------------------------

# ---- Option 1 heading ----
$ws.Range("A1").Value = "Option 1"

# ---- Option 1 table header (row 2) ----
$ws.Range("A2").Value = "SR NO"
$ws.Range("B2").Value = "ITEM DESCRIPTION"
$ws.Range("C2").Value = "QTY"
$ws.Range("D2").Value = "PRICE"
$ws.Range("E2").Value = "AMOUNT"
$hdr1 = $ws.Range("A2:E2")
$hdr1.Font.Bold = $true
$hdr1.HorizontalAlignment = -4108
$hdr1.VerticalAlignment = -4108
$hdr1.WrapText = $true
$hdr1.Borders.LineStyle = 1

# ---- Option 1 line items (rows 3-10) ----
$ws.Range("A3").Value = "1*"
$ws.Range("B3").Value = "Honeywell 2MP IP Bullet with inbuilt Audio"
$ws.Range("C3").Value = 16
$ws.Range("D3").Value = 4000
$ws.Range("E3").Formula = "=C3*D3"
$ws.Rows.Item(3).RowHeight = 22.2

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Honeywell NVR Professional Series 20CH"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 19900
$ws.Range("E4").Formula = "=C4*D4"
$ws.Rows.Item(4).RowHeight = 34.2

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "WD Purple Surveillance Hard Disk 4 TB"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 8600
$ws.Range("E5").Formula = "=C5*D5"
$ws.Rows.Item(5).RowHeight = 31.8

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "8 Port POE Giga Switch D link or Secue Eye or Similar"
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 8490
$ws.Range("E6").Formula = "=C6*D6"
$ws.Rows.Item(6).RowHeight = 26.4

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "RJ 45 Connector With Crimping etc complete"
$ws.Range("C7").Value = 35
$ws.Range("D7").Value = 150
$ws.Range("E7").Formula = "=C7*D7"
$ws.Rows.Item(7).RowHeight = 17.4

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Enclosure with mounting"
$ws.Range("C8").Value = 16
$ws.Range("D8").Value = 100
$ws.Range("E8").Formula = "=C8*D8"

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Supply and laying of cat 6 cables through PVC pipe / casing caping"
$ws.Range("C9").Value = 650
$ws.Range("D9").Value = 90
$ws.Range("E9").Formula = "=C9*D9"
$ws.Range("B9").WrapText = $true

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "INSTALLATION TESTING COMMISSIONING"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 5000
$ws.Range("E10").Formula = "=C10*D10"

$items1 = $ws.Range("A3:E10")
$items1.HorizontalAlignment = -4108
$items1.VerticalAlignment = -4108
$items1.WrapText = $true
$items1.Borders.LineStyle = 1

# ---- Option 1 TOTAL row (row 11, A:D merged) ----
$ws.Range("A11:D11").Merge()
$ws.Range("A11").Value = "TOTAL"
$ws.Range("E11").Formula = "=SUM(E3:E10)"
$tot1 = $ws.Range("A11:E11")
$tot1.Font.Bold = $true
$tot1.HorizontalAlignment = -4108
$tot1.VerticalAlignment = -4108
$tot1.Borders.LineStyle = 1

# ---- Notes below Option 1 table ----
$ws.Range("A13").Value = "In lieu of"
$ws.Range("A14").Value = "1) Honeywell 4 MP Bullet with inbuilt Audio 4 MM lens @ 6900+- GST"

$ws.Range("A16").Value = "If required extra Supply and laying of cat 6 cables through PVC pipe / casing caping to be chrarged as actualls @ 95/- per mtr + GST"

$ws.Range("A18").Value = "If required:-"
$ws.Range("A19").Value = "1)Display 19'' @7900 + GST"
$ws.Range("A20").Value = "2) Spike Board @ 600/- + GST"
$ws.Range("A21").Value = "3) Network Rack @ 2900/- + GST"
$ws.Range("A22").Value = "4) HDMI cable (3mtrs) @ 550/- + GST"

# ---------------------------------------------------------------------
# 2. Option 2 table (Dahua-based alternative)
# ---------------------------------------------------------------------
$ws.Range("A27").Value = "Option 2"
$ws.Range("A27").Font.Bold = $false

$ws.Range("A28").Value = "SR NO"
$ws.Range("B28").Value = "ITEM DESCRIPTION"
$ws.Range("C28").Value = "QTY"
$ws.Range("D28").Value = "PRICE"
$ws.Range("E28").Value = "AMOUNT"
$hdr2 = $ws.Range("A28:E28")
$hdr2.Font.Bold = $true
$hdr2.HorizontalAlignment = -4108
$hdr2.VerticalAlignment = -4108
$hdr2.WrapText = $true
$hdr2.Borders.LineStyle = 1

$ws.Range("A29").Value = 1
$ws.Range("B29").Value = "Dahus 2MP IP Bullet with inbuilt Audio"
$ws.Range("C29").Value = 16
$ws.Range("D29").Value = 4500
$ws.Range("E29").Formula = "=C29*D29"

$ws.Range("A30").Value = 3
$ws.Range("B30").Value = "Dahua NVR Professional Series 16CH"
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 13100
$ws.Range("E30").Formula = "=C30*D30"

$ws.Range("A31").Value = 3
$ws.Range("B31").Value = "WD Purple Surveillance Hard Disk 4` TB"
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 8600
$ws.Range("E31").Formula = "=C31*D31"

$ws.Range("A32").Value = 4
$ws.Range("B32").Value = "8 Port POE Giga Switch D link or Secue Eye or Similar"
$ws.Range("C32").Value = 4
$ws.Range("D32").Value = 8600
$ws.Range("E32").Formula = "=C32*D32"
$ws.Range("A32:E32").WrapText = $true

$ws.Range("A33").Value = 5
$ws.Range("B33").Value = "RJ 45 Connector With Crimping etc complete"
$ws.Range("C33").Value = 35
$ws.Range("D33").Value = 150
$ws.Range("E33").Formula = "=C33*D33"
$ws.Range("A33:E33").WrapText = $true

$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "Enclosure with mounting"
$ws.Range("C34").Value = 16
$ws.Range("D34").Value = 100
$ws.Range("E34").Formula = "=C34*D34"

$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "Supply and laying of cat 6 cables through PVC pipe / casing caping"
$ws.Range("C35").Value = 650
$ws.Range("D35").Value = 90
$ws.Range("E35").Formula = "=C35*D35"
$ws.Range("B35").WrapText = $true

$ws.Range("A36").Value = 8
$ws.Range("B36").Value = "INSTALLATION TESTING COMMISSIONING"
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = 5000
$ws.Range("E36").Formula = "=C36*D36"

$items2 = $ws.Range("A29:E36")
$items2.HorizontalAlignment = -4108
$items2.VerticalAlignment = -4108
$items2.Borders.LineStyle = 1

# ---- Option 2 TOTAL row (row 37, A:D merged) ----
$ws.Range("A37:D37").Merge()
$ws.Range("A37").Value = "TOTAL"
$ws.Range("E37").Formula = "=SUM(E29:E36)"
$tot2 = $ws.Range("A37:E37")
$tot2.Font.Bold = $true
$tot2.HorizontalAlignment = -4108
$tot2.VerticalAlignment = -4108
$tot2.Borders.LineStyle = 1

# ---- Notes below Option 2 table ----
$ws.Range("A38").Value = "In lieu of"
$ws.Range("A39").Value = "1) Dahua 4 MP Bullet camera @ 6500+- GST"

$ws.Range("A41").Value = "If required extra Supply and laying of cat 6 cables through PVC pipe / casing caping to be chrarged as actualls @ 95/- per mtr + GST"

$ws.Range("A43").Value = "If required:-"
$ws.Range("A44").Value = "1)Display 19'' @7900 + GST"
$ws.Range("A45").Value = "2) Spike Board @ 600/- + GST"
$ws.Range("A46").Value = "3) Network Rack @ 2900/- + GST"
$ws.Range("A47").Value = "4) HDMI cable (3mtrs) @ 550/- + GST"

# ---- GST note ----
$ws.Range("A57").Value = "levied"

# ---- trailing figures ----
$ws.Range("A60").Value = 75
$ws.Range("A61").Value = 20
$ws.Range("A62").Value = 5

# ---------------------------------------------------------------------
# 3. Make the new sheet the active / selected sheet (it was the target
#    of the upload, so Excel leaves focus on it, as reflected by the
#    workbook-level firstSheet/activeTab advancing by one tab).
# ---------------------------------------------------------------------
$ws.Select()
$ws.Range("H24").Select()

Write-Output "Gonsalves 204 sheet created"
